$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Row 6: Shipper Plantcode value changes from text "A001" to number 4552
$ws.Range("B6").Value = 4552

# Row 7: Destination Plantcode value changes from text "Z001" to text "Z999"
$ws.Range("B7").Value = "Z999"

# A10/A11 values updated
$ws.Range("A10").Value = 81730877
$ws.Range("A11").Value = 81730878

# Rows 12:15 - clear contents of columns A:B (C stays with its style/empty)
$ws.Range("A12:B15").ClearContents()

# Update selection to B8
$ws.Range("B8").Select()
